$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G10").Value = 3.4
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 2.15
$ws.Range("J10").Value = 4
$ws.Range("L10").Value = 2.88
$ws.Range("N10").Value = 9
$ws.Range("S10").Value = 1.5
$ws.Range("T10").Value = 2.5
$ws.Range("U10").Value = 1.95
$ws.Range("V10").Value = 1.8
$ws.Range("W10").Value = 9
$ws.Range("X10").Value = 17
$ws.Range("Z10").Value = 41
$ws.Range("AA10").Value = 29
$ws.Range("AD10").Value = 6.5
$ws.Range("AE10").Value = 17
$ws.Range("AG10").Value = 351
$ws.Range("AH10").Value = 6.5
$ws.Range("AI10").Value = 9.5
$ws.Range("AK10").Value = 19
$ws.Range("AL10").Value = 19
$ws.Range("AN10").Value = 5.5
$ws.Range("AO10").Value = 21
$ws.Range("AQ10").Value = 67
$ws.Range("AR10").Value = 101
$ws.Range("AS10").Value = 251
$ws.Range("AT10").Value = 2.5
$ws.Range("AU10").Value = 8.5
$ws.Range("AX10").Value = 4
$ws.Range("AY10").Value = 12
$ws.Range("AZ10").Value = 23
$ws.Range("G18").Value = 2.2
$ws.Range("J18").Value = 3
$ws.Range("Q18").Value = 2.4
$ws.Range("R18").Value = 1.53
$ws.Range("G19").Value = 1.67
$ws.Range("H19").Value = 3.5
$ws.Range("I19").Value = 5.75
$ws.Range("O19").Value = 1.36
$ws.Range("P19").Value = 3
$ws.Range("S19").Value = 1.44
$ws.Range("T19").Value = 2.63
$ws.Range("U19").Value = 2.1
$ws.Range("V19").Value = 1.67
$ws.Range("AL19").Value = 41
$ws.Range("AT19").Value = 2.63
$ws.Range("AY19").Value = 29
$ws.Range("S48").Value = 1.33
$ws.Range("G82").Value = 1.3
$ws.Range("J82").Value = 1.67
$ws.Range("K82").Value = 3.1
$ws.Range("L82").Value = 6.5
$ws.Range("AC82").Value = 34
$ws.Range("AF82").Value = 34
$ws.Range("AI82").Value = 51
$ws.Range("AJ82").Value = 23
$ws.Range("AK82").Value = 101
$ws.Range("AL82").Value = 51
$ws.Range("AR82").Value = 23
$ws.Range("AS82").Value = 51
$ws.Range("AW82").Value = 201
$ws.Range("BC82").Value = 101
